# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff run:
#   - new GUID-based file stem (f2a5d54c-cee0-42fd-9bb2-4ec3cd47035b)
#     replacing the previous one (4ae18501-3ac0-436c-a79b-804e5e220c5a)
#   - new content hash for the .xlf handoff payloads
#     (3bbf843518986509f7ecb940aefeabd8ff5ff5c7 replacing 1e8905f76d7eb58be35a978ae2c5f9b324b07861)
#   - refreshed handoff timestamps

$wb = $excel.ActiveWorkbook

$oldGuid = "4ae18501-3ac0-436c-a79b-804e5e220c5a"
$newGuid = "f2a5d54c-cee0-42fd-9bb2-4ec3cd47035b"
$oldHash = "1e8905f76d7eb58be35a978ae2c5f9b324b07861"
$newHash = "3bbf843518986509f7ecb940aefeabd8ff5ff5c7"

$newMdName = $newGuid + ".md"
$newZhXlfName = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newDeXlfName = $newGuid + "." + $newHash + ".de-de.xlf"

$newOverviewDate = "2016-03-24 19:20:30"
$newZhHandoffDatetime = "2016-03-24 19:20:26"
$newDeHandoffDatetime = "2016-03-24 19:20:30"

# External hyperlink targets stay the same as before the edit - only the
# displayed text (and the underlying cell values / shared strings) change.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/f9aad59642eb5a1af7b820168308c39faec35159/e2e/" + $oldGuid + ".md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4528ee908365cd3d011096b6d145f33b246c5f87/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $oldGuid + "." + $oldHash + ".zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f22df79f8751803c6cf7ab8928b9a15d9ae2c0be/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $oldGuid + "." + $oldHash + ".de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = $newZhHandoffDatetime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, "", "", $newZhXlfName)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = $newDeHandoffDatetime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, "", "", $newDeXlfName)
